$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1954397394136808
$ws.Range("C2").Value = 0.5504885993485342
$ws.Range("J2").Value = 0.02605863192182411
$ws.Range("P2").Value = 0.1302931596091205
$ws.Range("S2").Value = 0.09771986970684039

# Row 3
$ws.Range("C3").Value = 0.03409090909090909
$ws.Range("J3").Value = 0.03977272727272727
$ws.Range("P3").Value = 0.7045454545454546
$ws.Range("S3").Value = 0.2215909090909091

# Row 4
$ws.Range("J4").Value = 0.1132075471698113
$ws.Range("P4").Value = 0.6415094339622641
$ws.Range("S4").Value = 0.2452830188679245

# Row 6
$ws.Range("B6").Value = 0.0972972972972973
$ws.Range("D6").Value = 0.01621621621621622
$ws.Range("E6").Value = 0.005405405405405406
$ws.Range("F6").Value = 0.04324324324324325
$ws.Range("J6").Value = 0.2702702702702703
$ws.Range("O6").Value = 0.02162162162162162
$ws.Range("Q6").Value = 0.1567567567567568
$ws.Range("R6").Value = 0.04864864864864865
$ws.Range("S6").Value = 0.3405405405405406

# Row 7
$ws.Range("B7").Value = 0.1191709844559585
$ws.Range("D7").Value = 0.04663212435233161
$ws.Range("E7").Value = 0.005181347150259068
$ws.Range("F7").Value = 0.04663212435233161
$ws.Range("J7").Value = 0.1347150259067358
$ws.Range("O7").Value = 0.0155440414507772
$ws.Range("Q7").Value = 0.2020725388601036
$ws.Range("R7").Value = 0.1088082901554404
$ws.Range("S7").Value = 0.3212435233160622

# Row 8
$ws.Range("B8").Value = 0.09509202453987731
$ws.Range("D8").Value = 0.03680981595092025
$ws.Range("E8").Value = 0.003067484662576687
$ws.Range("F8").Value = 0.07668711656441718
$ws.Range("J8").Value = 0.147239263803681
$ws.Range("O8").Value = 0.01226993865030675
$ws.Range("Q8").Value = 0.1625766871165644
$ws.Range("R8").Value = 0.07668711656441718
$ws.Range("S8").Value = 0.3895705521472393

# Row 9
$ws.Range("B9").Value = 0.1420765027322404
$ws.Range("D9").Value = 0.01092896174863388
$ws.Range("F9").Value = 0.06010928961748634
$ws.Range("J9").Value = 0.1639344262295082
$ws.Range("O9").Value = 0.03278688524590164
$ws.Range("Q9").Value = 0.1420765027322404
$ws.Range("R9").Value = 0.0546448087431694
$ws.Range("S9").Value = 0.3934426229508197

# Row 10
$ws.Range("B10").Value = 0.1173139158576052
$ws.Range("D10").Value = 0.02346278317152103
$ws.Range("E10").Value = 0.0008090614886731392
$ws.Range("F10").Value = 0.06472491909385113
$ws.Range("J10").Value = 0.1440129449838188
$ws.Range("O10").Value = 0.01132686084142395
$ws.Range("Q10").Value = 0.244336569579288
$ws.Range("R10").Value = 0.07281553398058252
$ws.Range("S10").Value = 0.3211974110032362

# Row 11
$ws.Range("G11").Value = 0.1486486486486487
$ws.Range("J11").Value = 0.08445945945945946
$ws.Range("K11").Value = 0.1756756756756757
$ws.Range("L11").Value = 0.5743243243243243
$ws.Range("S11").Value = 0.01689189189189189

# Row 12
$ws.Range("G12").Value = 0.7167630057803468
$ws.Range("J12").Value = 0.2254335260115607
$ws.Range("K12").Value = 0.005780346820809248
$ws.Range("L12").Value = 0.03468208092485549
$ws.Range("S12").Value = 0.01734104046242774

# Row 13
$ws.Range("G13").Value = 0.6818181818181818
$ws.Range("J13").Value = 0.2727272727272727
$ws.Range("S13").Value = 0.04545454545454546

# Row 14
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333

# Row 15
$ws.Range("F15").Value = 0.03
$ws.Range("H15").Value = 0.13
$ws.Range("I15").Value = 0.075
$ws.Range("J15").Value = 0.395
$ws.Range("K15").Value = 0.065
$ws.Range("O15").Value = 0.095
$ws.Range("S15").Value = 0.21

# Row 16
$ws.Range("F16").Value = 0.0101010101010101
$ws.Range("H16").Value = 0.1717171717171717
$ws.Range("I16").Value = 0.06565656565656566
$ws.Range("J16").Value = 0.4242424242424243
$ws.Range("K16").Value = 0.1060606060606061
$ws.Range("M16").Value = 0.02525252525252525
$ws.Range("O16").Value = 0.0505050505050505
$ws.Range("S16").Value = 0.1464646464646465

# Row 17
$ws.Range("F17").Value = 0.01118568232662192
$ws.Range("H17").Value = 0.1610738255033557
$ws.Range("I17").Value = 0.0894854586129754
$ws.Range("J17").Value = 0.4049217002237137
$ws.Range("K17").Value = 0.116331096196868
$ws.Range("M17").Value = 0.01789709172259508
$ws.Range("O17").Value = 0.06935123042505593
$ws.Range("S17").Value = 0.1297539149888143

# Row 18
$ws.Range("F18").Value = 0.01948051948051948
$ws.Range("H18").Value = 0.1493506493506493
$ws.Range("I18").Value = 0.1428571428571428
$ws.Range("J18").Value = 0.4285714285714285
$ws.Range("K18").Value = 0.09740259740259741
$ws.Range("M18").Value = 0.01948051948051948
$ws.Range("O18").Value = 0.06493506493506493
$ws.Range("S18").Value = 0.07792207792207792

# Row 19
$ws.Range("F19").Value = 0.01596244131455399
$ws.Range("H19").Value = 0.1615023474178404
$ws.Range("I19").Value = 0.08544600938967137
$ws.Range("J19").Value = 0.3868544600938967
$ws.Range("K19").Value = 0.1342723004694836
$ws.Range("M19").Value = 0.02910798122065728
$ws.Range("N19").Value = 0.002816901408450704
$ws.Range("O19").Value = 0.06572769953051644
$ws.Range("S19").Value = 0.1183098591549296

Write-Output "Applied 111 changes"